$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rewrite the ODE-equation strings in column B for rows 12-24 so that every
# model id (column A) is paired with its corresponding, corrected equation
# text (missing "(t)" time-dependence markers have been added throughout).
# Column A (the model ids) and row 18 are left untouched - they were already
# correct.
# ---------------------------------------------------------------------------

$b12 = @'
odes = [
    sympy.Eq(S(t).diff(t), -beta * S(t) * I(t) / (1 + alpha*R(t))),
    sympy.Eq(I(t).diff(t), beta * S(t) * I(t) / (1 + alpha*R(t)) - gamma*I(t)),
    sympy.Eq(R(t).diff(t), gamma*I(t))
]
'@

$b13 = @'
odes = [
    sympy.Eq(S(t).diff(t), b - (beta_1*S(t)*P(t))/(1+alpha_1*P(t)) - (beta_2*S(t)*(I_A(t)+I_S(t)))/(1+alpha_2*(I_A(t)+I_S(t))) + psi*E(t) - mu*S(t)),
    sympy.Eq(E(t).diff(t), (beta_1*S(t)*P(t))/(1+alpha_1*P(t)) + (beta_2*S(t)*(I_A(t)+I_S(t)))/(1+alpha_2*(I_A(t)+I_S(t))) - psi*E(t) - mu*E(t) - omega*E(t)),
    sympy.Eq(I_A(t).diff(t), (1-delta)*omega*E(t) - (mu+sigma)*I_A(t) - gamma_A*I_A(t)),
    sympy.Eq(I_S(t).diff(t), delta*omega*E(t) - (mu+sigma)*I_S(t) - gamma_S*I_S(t)),
    sympy.Eq(R(t).diff(t), gamma_S*I_S(t) + gamma_A*I_A(t) - mu*R(t)),
    sympy.Eq(P(t).diff(t), eta_A*I_A(t) + eta_S*I_S(t) - mu_p*P(t))
]
'@

$b14 = @'
odes = [
    sympy.Eq(S(t).diff(t), -r_1*beta_1*I(t)*S(t)/N - r_2*beta_2*E(t)*S(t)/N),
    sympy.Eq(E(t).diff(t), r_1*beta*I(t)*S(t)/N - alpha*E(t) + r_2*beta_2*E(t)*S(t)/N),
    sympy.Eq(I(t).diff(t), alpha*E(t) - gamma*I(t)),
    sympy.Eq(R(t).diff(t), gamma*I(t))
]
'@

$b15 = @'
odes = [
    sympy.Eq(S(t).diff(t), -(beta*c + c*q*(1 - beta))*S(t)*(I(t) + theta*A(t)) + lambda_*S_q(t)),
    sympy.Eq(E(t).diff(t), beta*c*(1 - q)*S(t)*(I(t) + theta*A(t)) - sigma*E(t)),
    sympy.Eq(I(t).diff(t), sigma*rho*E(t) - (delta_I + alpha + gamma_I)*I(t)),
    sympy.Eq(A(t).diff(t), sigma*(1 - rho)*E(t) - gamma_A*A(t)),
    sympy.Eq(S_q(t).diff(t), (1 - beta)*c*q*S(t)*(I(t) + theta*A(t)) - lambda_*S_q(t)),
    sympy.Eq(E_q(t).diff(t), beta*c*q*S(t)*(I(t) + theta*A(t)) - delta_q*E_q(t)),
    sympy.Eq(H(t).diff(t), delta_I*I(t) + delta_q*E_q(t) - (alpha + gamma_H)*H(t)),
    sympy.Eq(R(t).diff(t), gamma_I*I(t) + gamma_A*A(t) + gamma_H*H(t))
]
'@

$b16 = @'
odes = [
    sympy.Eq(S(t).diff(t), -(beta*c(t) + c(t)*q*(1 - beta))*S(t)*(I(t) + theta*A(t)) + lambda_*S_q(t)),
    sympy.Eq(E(t).diff(t), beta*c(t)*(1 - q)*S(t)*(I(t) + theta*A(t)) - sigma*E(t)),
    sympy.Eq(I(t).diff(t), sigma*rho*E(t) - (delta_I(t) + alpha + gamma_I)*I(t)),
    sympy.Eq(A(t).diff(t), sigma*(1 - rho)*E(t) - gamma_A*A(t)),
    sympy.Eq(S_q(t).diff(t), (1 - beta)*c(t)*q*S(t)*(I(t) + theta*A(t)) - lambda_*S_q(t)),
    sympy.Eq(E_q(t).diff(t), beta*c(t)*q*S(t)*(I(t) + theta*A(t)) - delta_q*E_q(t)),
    sympy.Eq(H(t).diff(t), delta_I(t)*I(t) + delta_q*E_q(t) - (alpha + gamma_H)*H(t)),
    sympy.Eq(R(t).diff(t), gamma_I*I(t) + gamma_A*A(t) + gamma_H*H(t))
]
'@

$b17 = @'
odes = [
    sympy.Eq(S(t).diff(t), Lambda - mu*S(t) - beta*S(t)*I(t)/N),
    sympy.Eq(E(t).diff(t), beta*S(t)*I(t)/N - (mu + epsilon)*E(t)),
    sympy.Eq(I(t).diff(t), epsilon*E(t) - (gamma + mu + alpha)*I(t)),
    sympy.Eq(R(t).diff(t), gamma*I(t) - mu*R(t))
]
'@

$b19 = @'
odes = [
    sympy.Eq(S(t).diff(t), Lambda_s - (beta_s + rho_s*(1 - beta_s))*epsilon_s*S(t)*I(t)/N - delta*S(t) + m_s*S_q(t)),
    sympy.Eq(S_q(t).diff(t), (1 - beta_s)*epsilon_s*rho_s*S(t)*I(t)/N - (m_s + delta)*S_q(t)),
    sympy.Eq(A(t).diff(t), beta_s*(1 - rho_s)*epsilon_s*S(t)*I(t)/N - (gamma_a + xi_a + delta)*A(t)),
    sympy.Eq(I(t).diff(t), gamma_a*A(t) - (gamma_i + xi_i + delta)*I(t)),
    sympy.Eq(I_q(t).diff(t), beta_s*epsilon_s*rho_s*S(t)*I(t)/N + gamma_i*I(t) - (xi_q + delta)*I_q(t)),
    sympy.Eq(R(t).diff(t), xi_a*A(t) + xi_i*I(t) + xi_q*I_q(t) - delta*R(t))
]
'@

$b20 = @'
odes = [
    sympy.Eq(S(t).diff(t), -(1 - epsilon)*beta*S(t)*I(t)/N),
    sympy.Eq(E(t).diff(t), (1 - epsilon)*beta*S(t)*I(t)/N - sigma*E(t)),
    sympy.Eq(I(t).diff(t), sigma*E(t) - gamma*I(t)),
    sympy.Eq(R(t).diff(t), gamma*I(t))
]
'@

$b21 = @'
odes = [
    sympy.Eq(S(t).diff(t), -beta(t)*S(t)/N*I(t) + omega*R(t)),
    sympy.Eq(E(t).diff(t), beta(t)*S(t)/N*I(t) - sigma*E(t)),
    sympy.Eq(I(t).diff(t), sigma*E(t) - gamma*I(t)),
    sympy.Eq(R(t).diff(t), gamma*I(t) - omega*R(t))
]
'@

$b22 = @'
odes = [
    sympy.Eq(S_c(t).diff(t), m(t)*S_u(t) - (1 - m(t))*S_c(t)),
    sympy.Eq(S_u(t).diff(t), (1 - m(t))*S_c(t) - m(t)*S_u(t) - beta*(n*I_r(t) + I_u(t))*S_u(t) + theta*(1 - lambda_)*Q(t)),
    sympy.Eq(E(t).diff(t), (1 - sigma)*beta*(n*I_r(t) + I_u(t))*S_u(t) - mu*E(t)),
    sympy.Eq(I_r(t).diff(t), mu*f*E(t) + theta*lambda_*Q(t) - eta_r*I_r(t)),
    sympy.Eq(I_u(t).diff(t), mu*(1 - f)*E(t) - eta_u*I_u(t)),
    sympy.Eq(R(t).diff(t), eta_r*q*I_r(t) + eta_u*I_u(t)),
    sympy.Eq(Q(t).diff(t), sigma*beta*(n*I_r(t) + I_u(t))*S_u(t) - theta*Q(t))
]
'@

$b23 = @'
odes = [
    sympy.Eq(S(t).diff(t), -beta*S(t)/N*I(t)),
    sympy.Eq(E(t).diff(t), beta*S(t)/N*I(t) - omega*E(t)),
    sympy.Eq(I(t).diff(t), omega*E(t) - gamma*I(t)),
    sympy.Eq(R(t).diff(t), gamma*I(t))
]
'@

$b24 = @'
odes = [
    sympy.Eq(S(t).diff(t), -beta_c*(alpha*A(t) + I(t))*S(t)/(N_h - I_D(t))),
    sympy.Eq(E(t).diff(t), beta_c*(alpha*A(t) + I(t))*S(t)/(N_h - I_D(t)) - sigma*E(t)),
    sympy.Eq(A(t).diff(t), nu*sigma*E(t) - (theta + gamma_a)*A(t)),
    sympy.Eq(I(t).diff(t), (1 - nu)*sigma*E(t) - (psi + gamma_O + d_O)*I(t)),
    sympy.Eq(I_D(t).diff(t), theta*A(t) + psi*I(t) - (gamma_i + d_D)*I_D(t)),
    sympy.Eq(R(t).diff(t), gamma_i*I_D(t) + gamma_a*A(t) + gamma_O*I(t))
]
'@

$ws.Range("B12").Value2 = $b12
$ws.Range("B13").Value2 = $b13
$ws.Range("B14").Value2 = $b14
$ws.Range("B15").Value2 = $b15
$ws.Range("B16").Value2 = $b16
$ws.Range("B17").Value2 = $b17
$ws.Range("B19").Value2 = $b19
$ws.Range("B20").Value2 = $b20
$ws.Range("B21").Value2 = $b21
$ws.Range("B22").Value2 = $b22
$ws.Range("B23").Value2 = $b23
$ws.Range("B24").Value2 = $b24
# B18 already holds the correct, unchanged equation text - left untouched.

# ---------------------------------------------------------------------------
# Re-apply explicit row heights to match the corrected wrapped-text extents
# (only rows whose text actually changed length need a new height).
# ---------------------------------------------------------------------------
$ws.Rows.Item(12).RowHeight = 85
$ws.Rows.Item(13).RowHeight = 221
$ws.Rows.Item(14).RowHeight = 119
$ws.Rows.Item(15).RowHeight = 187
$ws.Rows.Item(16).RowHeight = 204
$ws.Rows.Item(17).RowHeight = 102
$ws.Rows.Item(18).RowHeight = 119
$ws.Rows.Item(19).RowHeight = 204
$ws.Rows.Item(20).RowHeight = 102
$ws.Rows.Item(21).RowHeight = 102
$ws.Rows.Item(22).RowHeight = 170
$ws.Rows.Item(23).RowHeight = 102
$ws.Rows.Item(24).RowHeight = 136

# ---------------------------------------------------------------------------
# Update the sheet view: the author scrolled/zoomed in and ended up with the
# selection on B22 while viewing near the bottom of the list at 150% zoom.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 150
$ws.Range("B22").Select()
